$d = $word.ActiveDocument

# The document has two paragraphs:
#   1) the heading "Sede Urbana P261 até P346" - must stay untouched
#   2) the long body paragraph listing the vertex-by-vertex description,
#      where every "P261".."P346" vertex label must be renumbered to
#      "P000".."P085" (new = old - 261, zero padded to 3 digits).
# Scope all Find/Replace calls to paragraph 2's Range so the heading is
# never touched.

$body = $d.Paragraphs.Item(2).Range

for ($old = 346; $old -ge 261; $old--) {
    $new = $old - 261
    $oldTag = "P" + $old
    $newTag = "P" + ("{0:D3}" -f $new)

    $rng = $body.Duplicate
    $rng.Find.Execute($oldTag, $false, $false, $false, $false, $false, $true, 1, $false, $newTag, 2)
}
